$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.411687
$ws.Range("H2").Value = 1.235061
$ws.Range("I2").Value = 0.00116706937075852
$ws.Range("J2").Value = 0.00116706937075852
$ws.Range("M2").Value = 1.116695
$ws.Range("N2").Value = 3.350085
$ws.Range("O2").Value = 0.008174214292497491
$ws.Range("P2").Value = 0.008174214292497492
$ws.Range("Q2").Value = 0.459728814465
$ws.Range("R2").Value = 4.137559330185
$ws.Range("S2").Value = [double]"9.53987513079035E-06"
$ws.Range("T2").Value = [double]"9.53987513079035E-06"
$ws.Range("G3").Value = 0.411687
$ws.Range("H3").Value = 1.235061
$ws.Range("I3").Value = 0.00116706937075852
$ws.Range("J3").Value = 0.00116706937075852
$ws.Range("O3").Value = 0.8193429796700005
$ws.Range("P3").Value = 0.8193429796700005
$ws.Range("Q3").Value = 46.08095202857999
$ws.Range("R3").Value = 414.72856825722
$ws.Range("S3").Value = 0.0009562300957188786
$ws.Range("T3").Value = 0.0009562300957188784
$ws.Range("G4").Value = 0.411687
$ws.Range("H4").Value = 1.235061
$ws.Range("I4").Value = 0.00116706937075852
$ws.Range("J4").Value = 0.00116706937075852
$ws.Range("O4").Value = 0.172482806037502
$ws.Range("P4").Value = 0.1724828060375021
$ws.Range("Q4").Value = 9.700665176834999
$ws.Range("R4").Value = 87.30598659151501
$ws.Range("S4").Value = 0.0002012993999088514
$ws.Range("T4").Value = 0.0002012993999088514
$ws.Range("G5").Value = 346.0613606666668
$ws.Range("I5").Value = 0.9810307695824356
$ws.Range("J5").Value = 0.9810307695824355
$ws.Range("M5").Value = 1.116695
$ws.Range("N5").Value = 3.350085
$ws.Range("O5").Value = 0.008174214292497491
$ws.Range("P5").Value = 0.008174214292497492
$ws.Range("Q5").Value = 386.4449911496635
$ws.Range("R5").Value = 3478.004920346971
$ws.Range("S5").Value = 0.008019155738100558
$ws.Range("T5").Value = 0.008019155738100558
$ws.Range("G6").Value = 346.0613606666668
$ws.Range("I6").Value = 0.9810307695824356
$ws.Range("J6").Value = 0.9810307695824355
$ws.Range("O6").Value = 0.8193429796700005
$ws.Range("P6").Value = 0.8193429796700005
$ws.Range("Q6").Value = 38735.34252921708
$ws.Range("S6").Value = 0.8038006738976264
$ws.Range("T6").Value = 0.8038006738976263
$ws.Range("G7").Value = 346.0613606666668
$ws.Range("I7").Value = 0.9810307695824356
$ws.Range("J7").Value = 0.9810307695824355
$ws.Range("O7").Value = 0.172482806037502
$ws.Range("P7").Value = 0.1724828060375021
$ws.Range("Q7").Value = 8154.314783967606
$ws.Range("R7").Value = 73388.83305570846
$ws.Range("S7").Value = 0.1692109399467086
$ws.Range("T7").Value = 0.1692109399467086
$ws.Range("I8").Value = 0.01780216104680593
$ws.Range("J8").Value = 0.01780216104680592
$ws.Range("M8").Value = 1.116695
$ws.Range("N8").Value = 3.350085
$ws.Range("O8").Value = 0.008174214292497491
$ws.Range("P8").Value = 0.008174214292497492
$ws.Range("Q8").Value = 7.012579198821666
$ws.Range("R8").Value = 63.11321278939499
$ws.Range("S8").Value = 0.0001455186792661431
$ws.Range("T8").Value = 0.0001455186792661431
$ws.Range("I9").Value = 0.01780216104680593
$ws.Range("J9").Value = 0.01780216104680592
$ws.Range("O9").Value = 0.8193429796700005
$ws.Range("P9").Value = 0.8193429796700005
$ws.Range("S9").Value = 0.01458607567665518
$ws.Range("T9").Value = 0.01458607567665518
$ws.Range("I10").Value = 0.01780216104680593
$ws.Range("J10").Value = 0.01780216104680592
$ws.Range("O10").Value = 0.172482806037502
$ws.Range("P10").Value = 0.1724828060375021
$ws.Range("S10").Value = 0.003070566690884601
$ws.Range("T10").Value = 0.003070566690884601
